$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Email column -------------------------------------------------
# Row 2 (Britney Spears): drop the stray space in the address.
# Rows 3 & 4 (Jennifer Lopez / Jennifer Aniston) were both erroneously
# showing Jennifer's own address - point them at the right person.
$ws.Range("C2").Value = "Britney@mail.com"
$ws.Range("C3").Value = "Lopez@mail.com"
$ws.Range("C4").Value = "Aniston@mail.com"

# --- Stamp the "Done on" column ------------------------------------------
$ws.Range("F2").Value = "8/18/2022  3:19:36 PM"
$ws.Range("F3").Value = "8/18/2022  3:20:00 PM"
$ws.Range("F4").Value = "8/18/2022  3:20:27 PM"

# --- Hyperlinks -----------------------------------------------------------
# The mailto hyperlinks keep pointing at the same addresses they always
# did; only the stale cached "display" text on the first one needs to go
# away. Remember the original "Hyperlink" cell formatting, recreate the
# three hyperlinks (which also drops the outdated display text) and then
# restore that formatting so the visible style doesn't change.
$ws.Range("C2").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Israel@mail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Jennifer@mail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Jennifer@mail.com")

$ws.Range("H1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Clear()

# --- Selection --------------------------------------------------------
# The workbook was last saved with C4 selected.
$ws.Range("C4").Select()
